$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the top for "Date and Time"
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 15:31:19.824000 to 2024-03-11 16:36:21.599000"

# 2. Insert a new row before "Idling time percentage" (currently row 35 after the shift above)
#    for "Cycle Count of battery"
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 114
